$wb = $excel.ActiveWorkbook

# Rename "baseline-forecasts" sheet to "external-forecasts"
$wsForecasts = $wb.Worksheets.Item("baseline-forecasts")
$wsForecasts.Name = "external-forecasts"

# Add trailing inflation variables to baseline-variables sheet
$wsVars = $wb.Worksheets.Item("baseline-variables")
$wsVars.Range("A23").Value = "dns1"
$wsVars.Range("A24").Value = "dns2"
$wsVars.Range("A25").Value = "dns3"

# Update selection to reflect where the user last clicked
$wsVars.Range("B23").Select()
